# Applies the COVID-19 Valais data-table edits described in the commit diff
# (rows 566-597 of sheet "Feuil1").
#
# Only the "hard-coded" input columns are written here:
#   C = Nb nouveaux cas positifs
#   E = Patients COVID-19 aux SI total (y.c. intubes)
#   F = Patients COVID-19 intubes
#   G = Patients COVID-19 hospitalises hors SI
#   L = Nb nouveaux deces a l'hopital
#   M = Nb nouveaux deces extra-hospitaliers
#
# Columns B, H, J, K hold live "IF(TODAY()>...)" shared formulas (ca="1") that
# depend on TODAY() and on the columns above, so their cached values update
# automatically on recalculation and need no explicit assignment.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns L and M are formatted as Text (numFmtId 49, "@"). Writing a number
# straight into ".Value" on such a cell gets stored as a text string (matches
# real Excel COM behaviour), which would wrongly flip the cell to t="s"/t="str".
# Round-tripping the number format to General for the write keeps the stored
# value numeric while leaving the cell's original (Text) display format intact
# once it's restored.
function Set-TextFormattedNumber {
    param($range, $value)
    $originalFormat = $range.NumberFormat
    $range.NumberFormat = "General"
    $range.Value = $value
    $range.NumberFormat = $originalFormat
}

# Row 566: F=6
$ws.Range("F566").Value = 6

# Row 567: F=3
$ws.Range("F567").Value = 3

# Row 568: E=8, F=4
$ws.Range("E568").Value = 8
$ws.Range("F568").Value = 4

# Row 569: E=8, F=2
$ws.Range("E569").Value = 8
$ws.Range("F569").Value = 2

# Row 570: E=8, F=1
$ws.Range("E570").Value = 8
$ws.Range("F570").Value = 1

# Row 571: E=7, F=1
$ws.Range("E571").Value = 7
$ws.Range("F571").Value = 1

# Row 572: E=7, F=1
$ws.Range("E572").Value = 7
$ws.Range("F572").Value = 1

# Row 573: E=7, F=1
$ws.Range("E573").Value = 7
$ws.Range("F573").Value = 1

# Row 574: E=5, F=1
$ws.Range("E574").Value = 5
$ws.Range("F574").Value = 1

# Row 575: E=5, F=1
$ws.Range("E575").Value = 5
$ws.Range("F575").Value = 1

# Row 576: E=4, F=1
$ws.Range("E576").Value = 4
$ws.Range("F576").Value = 1

# Row 577: E=4, F=1
$ws.Range("E577").Value = 4
$ws.Range("F577").Value = 1

# Row 578: E=4, F=1
$ws.Range("E578").Value = 4
$ws.Range("F578").Value = 1

# Row 579: E=3, F=1
$ws.Range("E579").Value = 3
$ws.Range("F579").Value = 1

# Row 580: E=3, F=1
$ws.Range("E580").Value = 3
$ws.Range("F580").Value = 1

# Row 581: E=3, F=1
$ws.Range("E581").Value = 3
$ws.Range("F581").Value = 1

# Row 582: E=2, F=0
$ws.Range("E582").Value = 2
$ws.Range("F582").Value = 0

# Row 583: C=46, E=2, F=0
$ws.Range("C583").Value = 46
$ws.Range("E583").Value = 2
$ws.Range("F583").Value = 0

# Row 584: E=1, F=0
$ws.Range("E584").Value = 1
$ws.Range("F584").Value = 0

# Row 585: E=1, F=0
$ws.Range("E585").Value = 1
$ws.Range("F585").Value = 0

# Row 586: E=1
$ws.Range("E586").Value = 1

# Row 587: E=1
$ws.Range("E587").Value = 1

# Row 588: E=1
$ws.Range("E588").Value = 1

# Row 589: E=1
$ws.Range("E589").Value = 1

# Row 590: C=55, E=1
$ws.Range("C590").Value = 55
$ws.Range("E590").Value = 1

# Row 591: C=35, E=1
$ws.Range("C591").Value = 35
$ws.Range("E591").Value = 1

# Row 592: C=23, E=1
$ws.Range("C592").Value = 23
$ws.Range("E592").Value = 1

# Row 593: E=2
$ws.Range("E593").Value = 2

# Row 594: C=60, E=2
$ws.Range("C594").Value = 60
$ws.Range("E594").Value = 2

# Row 595: C=55, E=2, F=1, G=7, L=0, M=0 (previously a fully blank placeholder row)
$ws.Range("C595").Value = 55
$ws.Range("E595").Value = 2
$ws.Range("F595").Value = 1
$ws.Range("G595").Value = 7
Set-TextFormattedNumber $ws.Range("L595") 0
Set-TextFormattedNumber $ws.Range("M595") 0

# Row 596: C=19, E=2, F=1, G=7, L=0, M=0 (previously a fully blank placeholder row)
$ws.Range("C596").Value = 19
$ws.Range("E596").Value = 2
$ws.Range("F596").Value = 1
$ws.Range("G596").Value = 7
Set-TextFormattedNumber $ws.Range("L596") 0
Set-TextFormattedNumber $ws.Range("M596") 0

# Row 597: C=2, E=2, F=1, G=6, L=0, M=0 (previously a fully blank placeholder row)
$ws.Range("C597").Value = 2
$ws.Range("E597").Value = 2
$ws.Range("F597").Value = 1
$ws.Range("G597").Value = 6
Set-TextFormattedNumber $ws.Range("L597") 0
Set-TextFormattedNumber $ws.Range("M597") 0
